$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking values
# (e.g. 0.998, 581.42) are not reinterpreted as numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.705.75"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.332.86"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "581.42"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "175.12"
$ws.Range("E6").Value = "  -5.93%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").Value = "3.329.04"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "0.177"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "45.59"
$ws.Range("E12").Value = "  -3.81%  "
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "666.00"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").Value = "3.866.16"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "67.787.00"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "3.326.91"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "17.48"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "0.890"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("D23").Value = "5.41"
$ws.Range("E23").Value = "  +6.26%  "
$ws.Range("D24").Value = "17.09"
$ws.Range("E24").Value = "  -5.06%  "
$ws.Range("D25").Value = "99.32"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "3.84"
$ws.Range("E26").Value = "  -6.01%  "
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -5.89%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "9.29"
$ws.Range("E28").Value = "  -5.59%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "33.79"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("D30").Value = "7.44"
$ws.Range("E30").Value = "  +8.63%  "
$ws.Range("D31").Value = "8.44"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").Value = "592.61"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("D33").Value = "10.96"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").Value = "3.710.17"
$ws.Range("E36").Value = "  -7.21%  "
$ws.Range("D37").Value = "56.86"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  -12.51%  "
$ws.Range("D39").Value = "0.132"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "33.29"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  -5.88%  "
$ws.Range("D42").Value = "3.11"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("D43").Value = "0.333"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").Value = "0.0₃0665"
$ws.Range("E44").Value = "  -5.78%  "
$ws.Range("D45").Value = "3.20"
$ws.Range("E45").Value = "  -6.59%  "
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").Value = "127.06"
$ws.Range("E51").Value = "  -3.21%  "

# Restore default (no explicit number format / style) on column D
$priceRange.Style = "Normal"

Write-Host "Updated cryptos list"
